# Generate Report for Handoff
# Insert a new "handoff" row (for 23db8880-cb21-4e34-b3f0-3b6f470c2e5d.md) as the
# new row 2 on every report sheet ("Overview", "zh-cn", "de-de"), pushing the
# existing b7c6459f-f78a-4417-a70c-2a960e35b2c3.md row down to row 3, then
# resize each sheet's table/hyperlinks to match.

$wb = $excel.ActiveWorkbook

$newGuid = "23db8880-cb21-4e34-b3f0-3b6f470c2e5d"
$newFile = "$newGuid.md"
$newPath = "e2e\$newFile"
$oldGuid = "b7c6459f-f78a-4417-a70c-2a960e35b2c3"
$oldFile = "$oldGuid.md"
$oldPath = "e2e\$oldFile"

$ghBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/010023c614bd69813348b0754b9f916ab4f990d8/"
$newUrl = $ghBase + $newPath
$oldUrl = $ghBase + $oldPath

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = $newPath
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-12 20:48:08"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("B2").Style = "Hyperlink"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $newUrl, "", "", $newPath)
$ws.Hyperlinks.Add($ws.Range("B3"), $oldUrl, "", "", $oldPath)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "$newGuid.07ba43b31c147c8fc63cb4e03a0905f1d88ba4ec.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-12 20:47:56"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""
$ws.Range("A2").Style = "Hyperlink"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $oldUrl, "", "", $oldFile)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "$newGuid.07ba43b31c147c8fc63cb4e03a0905f1d88ba4ec.de-de.xlf"
$ws.Range("H2").Value = "2016-08-12 20:48:08"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""
$ws.Range("A2").Style = "Hyperlink"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $oldUrl, "", "", $oldFile)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

Write-Output "Generated handoff report rows for $newFile"
